$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-10-08 10:42:41", "hatespeech", "def", 10, 0.9932370474835958),
    @("2023-10-08 10:42:41", "hatespeech", "def", 20, 0.9901418096316194),
    @("2023-10-08 10:42:41", "hatespeech", "def", 30, 0.9882337531908792),
    @("2023-10-08 10:42:41", "hatespeech", "def", 40, 0.9867270015936512)
)

$startRow = 23
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
